$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The previous hyperlinks (F2:F9) are anchored to fixed refs; inserting rows below
# will NOT relocate them, so drop them now and rebuild clean ones afterwards.
$ws.Range("A1:H9").Hyperlinks.Delete()

# Insert two new rows above the current row 2, shifting all data down by two rows.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Row 2
$ws.Range("A2").Value = '2025-12-21 18:25:16'
$ws.Range("B2").Value = '【TypeScript/Clasp必須】LINE WORKS連携・ファイル自動保存システムGAS開発'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5458419'
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5458419')
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("G2").Value = 193
$ws.Range("H2").Value = '🔥TypeScript ◆開発'

# Row 3
$ws.Range("A3").Value = '2025-12-21 18:25:16'
$ws.Range("B3").Value = '【急募】React/Supabaseで音楽権利マーケットプレイスMVP開発'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5458381'
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5458381')
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("G3").Value = 183
$ws.Range("H3").Value = '🔥React ◆開発'

# Row 4
$ws.Range("A4").Value = '2025-12-21 18:25:16'
$ws.Range("B4").Value = '【急募】ECサイトの自動購入Bot作成をお願いします。'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5458190'
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5458190')
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("G4").Value = 143
$ws.Range("H4").Value = '★bot ◇サイト'

# Row 5
$ws.Range("A5").Value = '2025-12-21 18:25:16'
$ws.Range("B5").Value = '【受注メールを元にECサイト自動仕入ツール】'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5458166'
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5458166')
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("G5").Value = 98
$ws.Range("H5").Value = '◆ツール ◇サイト'

# Row 6
$ws.Range("A6").Value = '2025-12-21 18:25:16'
$ws.Range("B6").Value = '【報酬計算の自動化】GASで自動計算させるプログラミング'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5458299'
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5458299')
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("G6").Value = 88
$ws.Range("H6").Value = '◆自動化'

# Row 7
$ws.Range("A7").Value = '2025-12-21 18:25:16'
$ws.Range("B7").Value = '【自動運転プロジェクト経験者募集】実証実験・開発を推進するプロジェクトマネージャー'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5431107'
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5431107')
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("G7").Value = 68
$ws.Range("H7").Value = '◆開発'

# Row 8
$ws.Range("A8").Value = '2025-12-21 18:25:16'
$ws.Range("B8").Value = '初回 【急募】ECサイトの要件定義や基本設計ができる方を募集(1人月、フルリモート可、2025年12月〜)'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5425629'
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5425629')
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("G8").Value = 45
$ws.Range("H8").Value = '◇サイト'

# Row 9
$ws.Range("A9").Value = '2025-12-21 18:25:16'
$ws.Range("B9").Value = '【急募】ManusアプリのGoogleStore登録代行を依頼したい'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5458330'
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5458330')
$ws.Range("F9").Style = "Hyperlink"
$ws.Range("G9").Value = 38
$ws.Range("H9").Value = '◇アプリ'

# Row 10
$ws.Range("A10").Value = '2025-12-21 18:25:16'
$ws.Range("B10").Value = '【急募】Notionでの社内向けダッシュボード作成依頼'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5458234'
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5458234')
$ws.Range("F10").Style = "Hyperlink"
$ws.Range("G10").Value = 18

# Row 11
$ws.Range("A11").Value = '2025-12-21 18:25:16'
$ws.Range("B11").Value = '【急募】グーグルワークスペースの設定をサポートしてくれる方'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '1,000 ~ 5,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5458288'
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5458288')
$ws.Range("F11").Style = "Hyperlink"
$ws.Range("G11").Value = 10

# Widen column H (17 chars) to match the new, longer "skill summary" text.
# ColumnWidth applies an internal padding of 5/6 of a character vs. the raw
# stored column width, so subtract that offset to land on exactly 17.
$ws.Columns.Item(8).ColumnWidth = 16.166666666666668

Write-Host "done"